$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for SVR parameters in columns K, L, M
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Add new SVR parameter values in row 2
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 20

# Row 13 previously held an empty styled placeholder cell; clear it out entirely
$ws.Range("A13").Clear()

# The duplicate "Normal" cell style that used to live alongside the built-in
# default got collapsed into the default style itself; drop the redundant
# explicit style from the cells that were only ever using that duplicate.
$ws.Range("B1").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("H1").Style = "Normal"
$ws.Range("H2").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").Style = "Normal"
$ws.Rows(7).ClearFormats()
$ws.Range("A8").Style = "Normal"
$ws.Range("A11").Style = "Normal"
$ws.Range("D14:G14").Style = "Normal"
$ws.Range("I14").Style = "Normal"
$ws.Range("D15:G15").Style = "Normal"
$ws.Range("I15").Style = "Normal"
$ws.Range("D16:G16").Style = "Normal"
$ws.Range("I16").Style = "Normal"

# Update the active selection to match the edited workbook
$ws.Range("K8").Select() | Out-Null
